# TC04_CCDI ... WebQuery sheet: update the "SamplesTab" SQL query (cell B4)
# so it also pulls participant_age_at_collection, sample_tumor_status and
# tumor_classification from df_sample, and surfaces them as three extra
# output columns ("Age at Sample Collection (days)", "Sample Tumor Status",
# "Sample Tumor Classification") before the final "Sample Diagnosis" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
WITH filtered_participants AS (
    SELECT DISTINCT prt.id, prt.participant_id
    FROM df_study std
    JOIN df_participant prt ON std.id = prt."study.id"
    JOIN df_diagnosis dgn ON prt.id = dgn."participant.id"
    WHERE 
        std.dbgap_accession = 'phs002371'
        AND prt.sex_at_birth = 'Not Reported'
        AND std.study_status = 'Completed'
),
-- Step 2: Get samples from those participants
filtered_samples AS (
    SELECT DISTINCT smp.id, smp.sample_id, smp.anatomic_site, smp."participant.id", smp.participant_age_at_collection, smp.sample_tumor_status,smp.tumor_classification
    FROM df_sample smp
    JOIN filtered_participants prt ON smp."participant.id" = prt.id
),
-- Step 3: Keep only samples that have at least one FASTQ file
samples_with_fastq AS (
    SELECT DISTINCT smp.*
    FROM filtered_samples smp
    JOIN df_sequencing_file seq ON seq."sample.id" = smp.id
    WHERE seq.file_type = 'bam'
)
-- Final output
SELECT 
    smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID",
    'phs002371' AS "Study ID",
    smp.anatomic_site AS "Sample Anatomic Site",
	COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '') AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '') AS "Sample Tumor Classification",
    NULL AS "Sample Diagnosis"
FROM 
    samples_with_fastq smp
JOIN filtered_participants prt ON smp."participant.id" = prt.id
ORDER BY 
    smp.sample_id ASC;
'@

$ws.Range("B4").Value = $newQuery

# The author's save also left the selection sitting on the edited cell.
$ws.Range("B4").Select() | Out-Null
